$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "A"
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "C"
$ws.Range("E1").Value = "D"
$ws.Range("F1").Value = "Message"
